# Generation & Import Working
# Replace the existing sample row (row 2) and append two more rows (3 & 4)
# of student records to the "student_information" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 -----------------------------------------------------------
$ws.Range("A2").Value = "Janet Scott"
$ws.Range("B2").Value = "deanna94@example.org"
$ws.Range("C2").Value = "557.987.3077x47839"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "6462397282"
$ws.Range("E2").Value = "Henrymouth"
$ws.Range("F2").Value = 3
$ws.Range("G2").Value = "Female"
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "2013-04-17"
$ws.Range("I2").Value = "1599 Amanda Plaza Suite 627, East Victoria, PW 17023"
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "2024-02-02"
$ws.Range("K2").Value = "South Dakota"
$ws.Range("L2").Value = "None"

# --- Row 3 (new) -------------------------------------------------------
$ws.Range("A3").Value = "Lynn Flores"
$ws.Range("B3").Value = "charles46@example.org"
$ws.Range("C3").Value = "+1-334-410-1697x5063"
$ws.Range("D3").Value = "+1-813-622-9373"
$ws.Range("E3").Value = "Jordanville"
$ws.Range("F3").Value = 12
$ws.Range("G3").Value = "Male"
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "2006-07-02"
$ws.Range("I3").Value = "62931 Chelsea Shore Suite 104, East George, IL 63484"
$ws.Range("J3").NumberFormat = "@"
$ws.Range("J3").Value = "2024-01-31"
$ws.Range("K3").Value = "North Carolina"
$ws.Range("L3").Value = "None"

# --- Row 4 (new) -------------------------------------------------------
$ws.Range("A4").Value = "Frank Castro"
$ws.Range("B4").Value = "hernandezcrystal@example.com"
$ws.Range("C4").Value = "(339)858-7240x786"
$ws.Range("D4").Value = "810-714-0034"
$ws.Range("E4").Value = "West Jeffrey"
$ws.Range("F4").Value = 7
$ws.Range("G4").Value = "Female"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "2014-03-23"
$ws.Range("I4").Value = "93634 James Lane, Hansenchester, MS 17361"
$ws.Range("J4").NumberFormat = "@"
$ws.Range("J4").Value = "2024-02-01"
$ws.Range("K4").Value = "Wyoming"
$ws.Range("L4").Value = "None"
